$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 4500
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 4500
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H98").Value = 8498.5
$ws.Range("I98").Value = 4148.8
$ws.Range("K98").Value = 4148.8
$ws.Range("M98").Value = -2650.8
$ws.Range("H105").Value = 67025
$ws.Range("J105").Value = 67025
$ws.Range("L105").Value = 67025
$ws.Range("N105").Value = -74013
$ws.Range("H122").Value = 8498.5
$ws.Range("I122").Value = 4148.8
$ws.Range("K122").Value = 12446.4
$ws.Range("M122").Value = -9996.400000000001
$ws.Range("H125").Value = 2500
$ws.Range("J125").Value = 2500
$ws.Range("L125").Value = 22500
$ws.Range("N125").Value = -27420
$ws.Range("H136").Value = 90000
$ws.Range("J136").Value = 90000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -100200
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 5559145
$ws.Range("I138").Value = 25003188
$ws.Range("J138").Value = 3704.0715
$ws.Range("K138").Value = 75009564
$ws.Range("L138").Value = 11112.2145
$ws.Range("M138").Value = -75004424
$ws.Range("N138").Value = -21392.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2577.7144
$ws.Range("I32").Value = 2577.7144
$ws.Range("K32").Value = 2577.7144
$ws.Range("M32").Value = -2290.7144
$ws.Range("H102").Value = 1078.1666
$ws.Range("I102").Value = 1078.1666
$ws.Range("K102").Value = 1078.1666
$ws.Range("M102").Value = 543.8334
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3999.6667
$ws.Range("J86").Value = 4999.5
$ws.Range("L86").Value = 4999.5
$ws.Range("N86").Value = -7245.5
$ws.Range("H89").Value = 3999.6667
$ws.Range("J89").Value = 4999.5
$ws.Range("L89").Value = 24997.5
$ws.Range("N89").Value = -36229.5
$ws.Range("H105").Value = 1943.5555
$ws.Range("I105").Value = 1784.7142
$ws.Range("K105").Value = 1784.7142
$ws.Range("M105").Value = -37.71419999999989
$ws.Range("H134").Value = 2937.8
$ws.Range("I134").Value = 3063.3333
$ws.Range("K134").Value = 9189.999899999999
$ws.Range("M134").Value = -6654.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 253.22223
$ws.Range("I7").Value = 314.57144
$ws.Range("K7").Value = 314.57144
$ws.Range("M7").Value = -201.57144
$ws.Range("H31").Value = 9299.75
$ws.Range("I31").Value = 7399.6665
$ws.Range("K31").Value = 7399.6665
$ws.Range("M31").Value = -7104.6665
$ws.Range("H32").Value = 3376.3333
$ws.Range("I32").Value = 1920.3334
$ws.Range("K32").Value = 1920.3334
$ws.Range("M32").Value = -1604.3334
$ws.Range("H34").Value = 9299.75
$ws.Range("I34").Value = 7399.6665
$ws.Range("K34").Value = 7399.6665
$ws.Range("M34").Value = -7197.6665
$ws.Range("H58").Value = 6256
$ws.Range("I58").Value = 6256
$ws.Range("K58").Value = 6256
$ws.Range("M58").Value = -6053
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 2256
$ws.Range("I99").Value = 2256
$ws.Range("K99").Value = 2256
$ws.Range("M99").Value = -758
$ws.Range("H106").Value = 26500
$ws.Range("J106").Value = 26500
$ws.Range("L106").Value = 26500
$ws.Range("N106").Value = -29024
$ws.Range("H126").Value = 2256
$ws.Range("I126").Value = 2256
$ws.Range("K126").Value = 6768
$ws.Range("M126").Value = -4298
$ws.Range("H132").Value = 12483
$ws.Range("I132").Value = 9949.5
$ws.Range("K132").Value = 29848.5
$ws.Range("M132").Value = -27318.5
$ws.Range("H136").Value = 6256
$ws.Range("I136").Value = 6256
$ws.Range("K136").Value = 18768
$ws.Range("M136").Value = -16218
$ws.Range("H141").Value = 441360.62
$ws.Range("J141").Value = 441360.62
$ws.Range("L141").Value = 441360.62
$ws.Range("N141").Value = -451720.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.916668
$ws.Range("I12").Value = 56.5
$ws.Range("K12").Value = 169.5
$ws.Range("M12").Value = 3.5
$ws.Range("H109").Value = 2449.6667
$ws.Range("I109").Value = 2449.6667
$ws.Range("K109").Value = 7349.000100000001
$ws.Range("M109").Value = -6309.000100000001
$ws.Range("H117").Value = 1855.25
$ws.Range("I117").Value = 5000
$ws.Range("K117").Value = 15000
$ws.Range("M117").Value = -11558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 22.933332
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 39.25
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 39.25
$ws.Range("M2").Value = 96
$ws.Range("N2").Value = -265.25
$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -31064
$ws.Range("H70").Value = 3499.5
$ws.Range("I70").Value = 3499.3333
$ws.Range("K70").Value = 3499.3333
$ws.Range("M70").Value = -3229.3333
$ws.Range("H73").Value = 3499.5
$ws.Range("I73").Value = 3499.3333
$ws.Range("K73").Value = 3499.3333
$ws.Range("M73").Value = -2563.3333
$ws.Range("H99").Value = 7900
$ws.Range("I99").Value = 7900
$ws.Range("K99").Value = 7900
$ws.Range("M99").Value = -5654
$ws.Range("H104").Value = 37835.5
$ws.Range("J104").Value = 37835.5
$ws.Range("L104").Value = 37835.5
$ws.Range("N104").Value = -44823.5
$ws.Range("H122").Value = 3000
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 8500
$ws.Range("I132").Value = 8500
$ws.Range("K132").Value = 25500
$ws.Range("M132").Value = -22970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5747.857
$ws.Range("I22").Value = 794.5
$ws.Range("J22").Value = 7729.2
$ws.Range("K22").Value = 794.5
$ws.Range("L22").Value = 7729.2
$ws.Range("M22").Value = -499.5
$ws.Range("N22").Value = -8319.200000000001
$ws.Range("H27").Value = 5747.857
$ws.Range("I27").Value = 794.5
$ws.Range("J27").Value = 7729.2
$ws.Range("K27").Value = 794.5
$ws.Range("L27").Value = 7729.2
$ws.Range("M27").Value = -687.5
$ws.Range("N27").Value = -7943.2
$ws.Range("H46").Value = 2305.2
$ws.Range("I46").Value = 2258.5
$ws.Range("J46").Value = 2375.25
$ws.Range("K46").Value = 2258.5
$ws.Range("L46").Value = 2375.25
$ws.Range("M46").Value = -2070.5
$ws.Range("N46").Value = -2751.25
$ws.Range("H68").Value = 2869.6667
$ws.Range("J68").Value = 2499
$ws.Range("L68").Value = 2499
$ws.Range("N68").Value = -3997
$ws.Range("H71").Value = 2869.6667
$ws.Range("J71").Value = 2499
$ws.Range("L71").Value = 12495
$ws.Range("N71").Value = -19983
$ws.Range("H93").Value = 1593
$ws.Range("I93").Value = 1450.25
$ws.Range("K93").Value = 1450.25
$ws.Range("M93").Value = -202.25
$ws.Range("H100").Value = 2633.3333
$ws.Range("I100").Value = 2633.3333
$ws.Range("K100").Value = 2633.3333
$ws.Range("M100").Value = -2092.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1543
$ws.Range("I122").Value = 1428.9166
$ws.Range("J122").Value = 1999.3334
$ws.Range("K122").Value = 4286.7498
$ws.Range("L122").Value = 5998.0002
$ws.Range("M122").Value = -1836.7498
$ws.Range("N122").Value = -10898.0002
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H133").Value = 69998
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
